$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "290.43"
$ws.Range("E2").Value = "-6.08%"
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "39.82"
$ws.Range("E3").Value = "-3.00%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").Value = "-3.42%"
$rng.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "-4.08%"
$cell.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "4.285"
$ws.Range("E6").Value = "-0.37%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "1.549"
$ws.Range("E7").Value = "-11.62%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9150"
$ws.Range("E8").Value = "-0.05%"
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1189"
$ws.Range("E9").Value = "-5.78%"
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.1736"
$ws.Range("E10").Value = "-4.70%"
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.08691"
$ws.Range("E11").Value = "-4.51%"
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.04184"
$ws.Range("E12").Value = "0.57%"
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").Value = "0.02%"
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001277"
$ws.Range("E14").Value = "-0.26%"
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.005883"
$ws.Range("E15").Value = "-0.21%"
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "3.390"
$ws.Range("E16").Value = "1.08%"
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "0.3295"
$ws.Range("E18").Value = "-0.74%"
$rng.Style = "Normal"

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "7.550"
$ws.Range("E19").Value = "1.83%"
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = "0.1352"
$ws.Range("E20").Value = "-0.21%"
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "0.2736"
$ws.Range("E21").Value = "0.40%"
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "0.03835"
$ws.Range("E22").Value = "-4.39%"
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.001276"
$ws.Range("E23").Value = "0.59%"
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.003875"
$ws.Range("E24").Value = "-5.44%"
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.0001283"
$ws.Range("E25").Value = "0.80%"
$rng.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0003732"
$cell.Style = "Normal"

$rng = $ws.Range("D38:E38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = "0.02321"
$ws.Range("E38").Value = "-7.94%"
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.05019"
$ws.Range("E39").Value = "-5.14%"
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.007693"
$ws.Range("E40").Value = "-2.07%"
$rng.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "172.53%"
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "-3.16%"
$cell.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.007410"
$ws.Range("E43").Value = "11.25%"
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.007672"
$ws.Range("E44").Value = "-5.72%"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.3163"
$ws.Range("E45").Value = "3.02%"
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006520"
$ws.Range("E46").Value = "-4.01%"
$rng.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "0.01%"
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "9.55%"
$cell.Style = "Normal"

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = "0.004211"
$ws.Range("E49").Value = "35.68%"
$rng.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "0.01%"
$cell.Style = "Normal"

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.01%"
$rng.Style = "Normal"
